$wb = $excel.ActiveWorkbook

# --- Withdraw History: replace Day/Month/Year/Place columns with Date/Location-Place ---
$ws2 = $wb.Worksheets.Item("Withdraw History")
$ws2.Range("C1").Value = "Date"
$ws2.Range("D1").Value = "Location/Place"
$ws2.Range("E1:F1").ClearContents()

# --- Deposit History: same remapping as Withdraw History ---
$ws3 = $wb.Worksheets.Item("Deposit History")
$ws3.Range("C1").Value = "Date"
$ws3.Range("D1").Value = "Location/Place"
$ws3.Range("E1:F1").ClearContents()

# --- Transfer History: Date then Person ---
$ws4 = $wb.Worksheets.Item("Transfer History")
$ws4.Range("C1").Value = "Date"
$ws4.Range("D1").Value = "Person"
$ws4.Range("E1:F1").ClearContents()

# --- Absolute History: Date then Location/Place/Person ---
$ws5 = $wb.Worksheets.Item("Absolute History")
$ws5.Range("C1").Value = "Date"
$ws5.Range("D1").Value = "Location/Place/Person"
$ws5.Range("E1:F1").ClearContents()

# --- Make "Absolute History" the active tab (moves tabSelected + sets activeTab=4) ---
$ws5.Activate()
